$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 66.63
$ws.Range("H2").Value = 66.63
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 70
$ws.Range("F3").Value = 65.85
$ws.Range("G3").Value = 70
$ws.Range("H3").Value = 65.85
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 66.67
$ws.Range("D4").Value = 80
$ws.Range("E4").Value = 156.25
$ws.Range("F4").Value = 44.48
$ws.Range("G4").Value = 93.75
$ws.Range("H4").Value = 40.24
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 77.78
$ws.Range("D5").Value = 80.56
$ws.Range("E5").Value = 104
$ws.Range("F5").Value = 22.88
$ws.Range("G5").Value = 84
$ws.Range("H5").Value = 19.88
$ws.Range("B6").Value = 10
$ws.Range("C6").Value = 90
$ws.Range("D6").Value = 95.73
$ws.Range("E6").Value = 85
$ws.Range("F6").Value = 32.72
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = 32.69
$ws.Range("B7").Value = 13
$ws.Range("C7").Value = 84.62
$ws.Range("D7").Value = 77.84
$ws.Range("E7").Value = 114
$ws.Range("F7").Value = 56.52
$ws.Range("G7").Value = 96.67
$ws.Range("H7").Value = 30.71
$ws.Range("B8").Value = 14
$ws.Range("C8").Value = 85.71
$ws.Range("D8").Value = 78.15
$ws.Range("E8").Value = 120.49
$ws.Range("F8").Value = 87.1
$ws.Range("G8").Value = 99.18
$ws.Range("H8").Value = 35.91
$ws.Range("C9").Value = 20
$ws.Range("D9").Value = 77.52
$ws.Range("E9").Value = 131
$ws.Range("F9").Value = 39.19
$ws.Range("G9").Value = 100
$ws.Range("H9").Value = 16.54
$ws.Range("B10").Value = 2
$ws.Range("E10").Value = 85
$ws.Range("F10").Value = 73.97
$ws.Range("G10").Value = 85
$ws.Range("H10").Value = 73.97
$ws.Range("B11").Value = 20
$ws.Range("C11").Value = 70
$ws.Range("D11").Value = 77.78
$ws.Range("E11").Value = 180
$ws.Range("F11").Value = 183.94
$ws.Range("G11").Value = 100
$ws.Range("H11").Value = 88.94
$ws.Range("B12").Value = 2
$ws.Range("D12").Value = 77.52
$ws.Range("E12").Value = 70
$ws.Range("F12").Value = 22.88
$ws.Range("G12").Value = 70
$ws.Range("H12").Value = 16.54
$ws.Range("B13").Value = 20
$ws.Range("E13").Value = 180
$ws.Range("F13").Value = 183.94
$ws.Range("H13").Value = 88.94
$ws.Range("B14").Value = 8.9
$ws.Range("C14").Value = 79.47800000000001
$ws.Range("D14").Value = 86.758
$ws.Range("E14").Value = 114.574
$ws.Range("F14").Value = 67.328
$ws.Range("G14").Value = 90.35999999999999
$ws.Range("H14").Value = 47.136
$ws.Range("B15").Value = 5.466056876558985
$ws.Range("C15").Value = 24.11509660864838
$ws.Range("D15").Value = 10.59458855579898
$ws.Range("E15").Value = 33.9654854881311
$ws.Range("F15").Value = 45.50947124866794
$ws.Range("G15").Value = 11.19289357881449
$ws.Range("H15").Value = 24.78563840434841
$ws.Range("B16").Value = 30.36698264754992
$ws.Range("C16").Value = 30.14387076081047
$ws.Range("D16").Value = 47.12895265035132
$ws.Range("E16").Value = 30.87771408011918
$ws.Range("F16").Value = 28.25622205927477
$ws.Range("G16").Value = 37.30964526271497
$ws.Range("H16").Value = 34.23430718832653
